$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (date serial 45016 = 2023-03-31) is added
# right before the current row 272, pushing every following data row down
# by one (old row 272 -> new row 273, ..., old row 409 -> new row 410).
$ws.Rows(272).Insert()

# Populate the newly inserted row 272 with the new weekly record. The
# non-date fields mirror the values that row 272 held prior to the shift
# (now living in row 273), per the dataset's weekly-refresh pattern.
$ws.Range("A272").Value = 9
$ws.Range("B272").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C272").Value = "Metropolitana"
$ws.Range("D272").Value = 45016
$ws.Range("E272").Value = 13
$ws.Range("F272").Value = 300000001
$ws.Range("G272").Value = "Rabanito"
$ws.Range("H272").Value = "Sin especificar"
$ws.Range("I272").Value = "Primera"
$ws.Range("J272").Value = 7000
$ws.Range("K272").Value = 3000
$ws.Range("L272").Value = 3000
$ws.Range("M272").Value = 3000
$ws.Range("N272").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O272").Value = "Provincia de Chacabuco"
$ws.Range("P272").Value = 30
$ws.Range("Q272").Value = 100
$ws.Range("R272").Value = "Hortaliza"
